$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.132418990135193
$ws.Range("B1").Value = 2.286110401153564
$ws.Range("C1").Value = 11.05879497528076
$ws.Range("D1").Value = 2.086678266525269
$ws.Range("E1").Value = 1.278679966926575
